# Insert a new "StatQuery" column between the existing "query" column (A)
# and the "dbExcel" column (old B, now shifts to C), giving the layout:
#   A: query        B: StatQuery   C: dbExcel           D: WebExcel
#   A2: <long match> B2: <stat query> C2: Neo4jData.xlsx D2: WebData.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; this shifts the old B ("dbExcel"/Neo4jData.xlsx)
# to C and old C ("WebExcel"/WebData.xlsx) to D, while preserving their
# column widths/content automatically.
$ws.Columns("B").Insert()

# The insert copies the neighboring (wrap-text) style into the new column;
# reset it back to the default "Normal" style, matching the unstyled
# original B/C columns.
$ws.Range("B1:B2").Style = "Normal"

# Header for the newly inserted column.
$ws.Range("B1").Value2 = "StatQuery"

# New stat-bar query text for row 2 of the new column.
$statQuery = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Bouvier des Flandres']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@
$ws.Range("B2").Value2 = $statQuery

# Widen the new column so the long query text is fully visible (matches the
# author's manual "best fit" resize of column B to ~255.6 characters).
$ws.Columns("B").ColumnWidth = 254.8333333

# Update the active selection to B2 (the cell the author last edited/selected).
$null = $ws.Range("B2").Select()
